$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A95").Value = "GRT-USD"
